$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.110.25'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.548.22'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.03%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.86'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.36%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.546.20'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.85%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.105'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.56'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.13%  '
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.352'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.57'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.003.36'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.93%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.061.91'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.29%  '
$ws.Range('E17').Value = '  +1.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.555.47'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.42'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.56%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '335.30'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.32'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.54%  '
$ws.Range('E22').Value = '  -0.69%  '
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.17'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.63'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +8.97%  '
$ws.Range('E26').Value = '  -1.93%  '
$ws.Range('E27').Value = '  +8.36%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.45'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.36'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0819'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.57%  '
$ws.Range('E32').Value = '  +0.63%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '176.05'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.57'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '412.33'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +11.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.400'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.58%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.95'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.40'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.40%  '
$ws.Range('E40').Value = '  +3.82%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.33'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '152.73'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.78'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.07'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.605'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.63%  '
$ws.Range('E47').Value = '  +0.43%  '
$ws.Range('E48').Value = '  +1.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0238'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.36'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.12%  '
$ws.Range('E51').Value = '  +1.09%  '
